# Chabbi_2009: Updated USDA soil type
#
# Adds a new "pro_usda_soil_order" field to the "profile" sheet (inserted
# right after pro_MAP / before pro_soil_taxon), records the USDA soil
# order + soil taxon / taxonomic system for the Chabbi_2009 profile row,
# and adds the corresponding controlled-vocabulary column (with the 12
# USDA soil orders) on the "controlled vocabulary" sheet.

$wb = $excel.ActiveWorkbook

# --- profile sheet: insert new column N (pro_usda_soil_order) ---
$wsProfile = $wb.Worksheets.Item("profile")
$wsProfile.Columns.Item(14).Insert()

$wsProfile.Cells.Item(4, 15).Value = "Dystric Cambisol"
$wsProfile.Cells.Item(1, 14).Value = "pro_usda_soil_order"

# --- controlled vocabulary sheet: insert new column E (pro_usda_soil_order list) ---
$wsCV = $wb.Worksheets.Item("controlled vocabulary")
$wsCV.Columns.Item(5).Insert()

$wsCV.Cells.Item(4, 5).Value = "Alfisols"
$wsCV.Cells.Item(5, 5).Value = "Andisols"
$wsCV.Cells.Item(6, 5).Value = "Aridisols"
$wsCV.Cells.Item(7, 5).Value = "Entisols"
$wsCV.Cells.Item(8, 5).Value = "Gelisols"
$wsCV.Cells.Item(9, 5).Value = "Histosols"

$wsProfile.Cells.Item(4, 14).Value = "Inceptisols"

$wsCV.Cells.Item(11, 5).Value = "Mollisols"
$wsCV.Cells.Item(12, 5).Value = "Oxisols"
$wsCV.Cells.Item(13, 5).Value = "Spodosols"
$wsCV.Cells.Item(14, 5).Value = "Ultisols"
$wsCV.Cells.Item(15, 5).Value = "Vertisols"

$wsCV.Cells.Item(10, 5).Value = "Inceptisols"

$wsProfile.Cells.Item(4, 17).Value = "WRB"

$wsCV.Cells.Item(2, 5).Value = "pro_usda_soil_order"
